# Updates the "cryptos" price/volume table with freshly scraped values.
# (Updated cryptos list on Fri Aug 18 07:00:02 UTC 2023 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold numeric-looking text (e.g. "1.005", "26.630.95")
# that must stay plain text, exactly like the original inline strings, instead of
# being auto-converted to numbers by Excel. Temporarily force a text format over
# the whole data range, write the new values, then restore the default "Normal"
# style so the cells end up unformatted again, just like before the edit.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.630.95"
$ws.Range("E2").Value = "  -7.40%  "
$ws.Range("D3").Value = "1.697.33"
$ws.Range("E3").Value = "  -6.07%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "219.41"
$ws.Range("E5").Value = "  -5.57%  "
$ws.Range("D6").Value = "0.5140"
$ws.Range("E6").Value = "  -12.99%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.2660"
$ws.Range("E8").Value = "  -4.29%  "
$ws.Range("D9").Value = "22.17"
$ws.Range("E9").Value = "  -4.96%  "
$ws.Range("D10").Value = "0.06255"
$ws.Range("E10").Value = "  -8.52%  "
$ws.Range("D11").Value = "0.07329"
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("D12").Value = "1.698.69"
$ws.Range("E12").Value = "  -6.10%  "
$ws.Range("D13").Value = "4.521"
$ws.Range("E13").Value = "  -5.10%  "
$ws.Range("D14").Value = "0.5848"
$ws.Range("E14").Value = "  -6.22%  "
$ws.Range("D15").Value = "1.927.85"
$ws.Range("E15").Value = "  -6.06%  "
$ws.Range("D16").Value = "0.000008427"
$ws.Range("E16").Value = "  -9.25%  "
$ws.Range("D17").Value = "65.51"
$ws.Range("E17").Value = "  -13.55%  "
$ws.Range("D18").Value = "26.674.63"
$ws.Range("E18").Value = "  -7.11%  "
$ws.Range("D19").Value = "5.049"
$ws.Range("E19").Value = "  -7.85%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "10.90"
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("D22").Value = "187.29"
$ws.Range("E22").Value = "  -11.50%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "145.21"
$ws.Range("E25").Value = "  -5.79%  "
$ws.Range("D26").Value = "7.614"
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("E27").Value = "  -9.43%  "
$ws.Range("D28").Value = "15.81"
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("D29").Value = "1.309"
$ws.Range("E29").Value = "  -8.93%  "
$ws.Range("D30").Value = "0.05729"
$ws.Range("E30").Value = "  -7.55%  "
$ws.Range("D31").Value = "1.333"
$ws.Range("E31").Value = "  -6.55%  "
$ws.Range("D32").Value = "3.528"
$ws.Range("E32").Value = "  -6.19%  "
$ws.Range("D33").Value = "3.510"
$ws.Range("E33").Value = "  -7.37%  "
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("D36").Value = "0.6031"
$ws.Range("E36").Value = "  -6.32%  "
$ws.Range("D37").Value = "2.377"
$ws.Range("E37").Value = "  -4.76%  "
$ws.Range("D38").Value = "2.683"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").Value = "1.097.45"
$ws.Range("E39").Value = "  -4.35%  "
$ws.Range("D40").Value = "0.01602"
$ws.Range("E40").Value = "  -5.86%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "5.904"
$ws.Range("E41").Value = "  -10.16%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8644"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "98.78"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").Value = "1.855.81"
$ws.Range("E45").Value = "  -5.38%  "
$ws.Range("D46").Value = "0.00000000108"
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("D47").Value = "56.85"
$ws.Range("E47").Value = "  -6.16%  "
$ws.Range("D48").Value = "8.157"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "0.05245"
$ws.Range("E50").Value = "  -4.20%  "
$ws.Range("E51").Value = "  -3.55%  "

$dataRange.Style = "Normal"
